$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update summary header values ---
# VALOR MORA total
$ws.Range("E11").Value = 27041
# Cant. Trabajadores
$ws.Range("C13").Value = 3
# Cant. Periodos
$ws.Range("F13").Value = 2

# --- Reorder the three worker detail rows (16-18) ---
# New order: HAILER DE LA ROSA MONTIEL, ARAMIS MENDOZA VALDEZ, YEISER EDUARDO AVILA HEREDIA
$ws.Range("B16").Value = "CC"
$ws.Range("C16").Value = "73213306"
$ws.Range("D16").Value = "HAILER DE LA ROSA MONTIEL"
$ws.Range("E16").Value = "1807"
$ws.Range("F16").Value = 9375
$ws.Range("G16").Value = 781242

$ws.Range("B17").Value = "CC"
$ws.Range("C17").Value = "1007855319"
$ws.Range("D17").Value = "ARAMIS MENDOZA VALDEZ"
$ws.Range("E17").Value = "1901"
$ws.Range("F17").Value = 16562
$ws.Range("G17").Value = 828116

$ws.Range("B18").Value = "CC"
$ws.Range("C18").Value = "1047422179"
$ws.Range("D18").Value = "YEISER EDUARDO AVILA HEREDIA"
$ws.Range("E18").Value = "1901"
$ws.Range("F18").Value = 1104
$ws.Range("G18").Value = 828116

# Row 18 becomes the new bottom row of the table, so it must carry the
# "closing border" formatting that currently belongs to row 28 (the last
# row of the table before the old records are removed).
$ws.Range("B28:J28").Copy()
$ws.Range("B18:J18").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# --- Remove the old MONICA TAMAYO CASTAÑO rows (19-28) entirely ---
$ws.Range("A19:A28").EntireRow.Delete()
